$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 text with the new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$rngA1 = $wsHoja1.Range("A1")
$oldText = $rngA1.Value()
$newText = $oldText.Replace(
    "1000 Bs = 8.93 = 37200.08 pesos",
    "1000 Bs = 8.93 = 37503.86 pesos"
)
$newText = $newText.Replace(
    "37200.08 pesos = 8.87 = 948.41 Bs",
    "37503.86 pesos = 8.94 = 974.94 Bs"
)
$rngA1.Value = $newText

# --- Update numeric rate cells on the "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("O10").Value = 4200.02
$wsTasas.Range("O12").Value = 109
